# Daily attendance processing - 2025-12-09 11:50:02
# Normalize the "Recorded By" (column G) entries: when the literal "System"
# entry is listed first among a comma-separated list of recorders, move it
# so it swaps places with the entry immediately after it (i.e. it is no
# longer reported first in the list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ($val -ne "") {
        $parts = $val -split ", "
        if ($parts.Length -ge 2 -and $parts[0] -eq "System") {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $cell.Value = $parts -join ", "
        }
    }
}
